# User View -> Bulk User - users data array
#
# - Rebuilds the "Roles" sheet as a Role Code / Platform - Role lookup table
#   (4 roles: C150/AES - Admin AES, C151/TBS - Admin TBS, C152/TBS - Tenant,
#   C153/AES - Admin), replacing the old 2-row numeric-code table.
# - Adds workbook-level named ranges "users" (UserInfo!A2:A1048576) and
#   "roleCode" (Roles!A2:A1048576).
# - Adds list data-validation dropdowns on the "UserRoles" sheet driven by
#   those named ranges.
# - Restores per-sheet selection/active-cell state.

$wb = $excel.ActiveWorkbook

$wsUserInfo  = $wb.Worksheets.Item("UserInfo")
$wsUserRoles = $wb.Worksheets.Item("UserRoles")
$wsRoles     = $wb.Worksheets.Item("Roles")

# --- Roles sheet: rebuild as a Role Code / Platform - Role lookup table ---
$wsRoles.Range("A1").Value = "Role Code"
$wsRoles.Range("B1").Value = "Platform - Role"

$wsRoles.Range("A2").Value = "C150"
$wsRoles.Range("B2").Value = "AES - Admin AES"

$wsRoles.Range("A3").Value = "C151"
$wsRoles.Range("B3").Value = "TBS - Admin TBS"

$wsRoles.Range("A4").Value = "C152"
$wsRoles.Range("B4").Value = "TBS - Tenant"

$wsRoles.Range("A5").Value = "C153"
$wsRoles.Range("B5").Value = "AES - Admin"

$wsRoles.Columns.Item(1).ColumnWidth = 9.166666666666666
$wsRoles.Columns.Item(2).ColumnWidth = 22.666666666666668

# --- Workbook-level defined names used by the validation dropdowns ---
$wb.Names.Add("roleCode", "=Roles!`$A`$2:`$A`$1048576")
$wb.Names.Add("users", "=UserInfo!`$A`$2:`$A`$1048576")

# --- UserRoles sheet: list-style data validation dropdowns ---
$wsUserRoles.Range("A2:A1048576").Validation.Add(3, 1, 1, "=users")
$wsUserRoles.Range("B2:B1048576").Validation.Add(3, 1, 1, "=roleCode")

# --- Selections / active cells left behind on each sheet ---
$null = $wsUserInfo.Range("H5").Select()
$null = $wsUserRoles.Range("I19").Select()
$null = $wsRoles.Range("J21").Select()

# Roles is the tab that ends up selected/active.
$wsRoles.Activate()

$wb.Save()
